# Croatia 1NL - base update (28-05-2024 20:56)
# The underlying source data got re-sorted/re-paired for several fixtures that
# share the same match date. Column A (the running index) and the row
# position stay put; everything from column B through AD moves between the
# affected rows. This script captures each affected row's B:AD block first
# (so the reads aren't clobbered by earlier writes) and then redistributes
# them according to the new pairing/rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colRange = "B{0}:AD{0}"

# ---- snapshot every row that moves, before writing anything back ----
$row11 = $ws.Range(($colRange -f 11)).Value2
$row12 = $ws.Range(($colRange -f 12)).Value2

$row15 = $ws.Range(($colRange -f 15)).Value2
$row16 = $ws.Range(($colRange -f 16)).Value2
$row17 = $ws.Range(($colRange -f 17)).Value2

$row21 = $ws.Range(($colRange -f 21)).Value2
$row22 = $ws.Range(($colRange -f 22)).Value2

$row101 = $ws.Range(($colRange -f 101)).Value2
$row102 = $ws.Range(($colRange -f 102)).Value2

$row164 = $ws.Range(($colRange -f 164)).Value2
$row165 = $ws.Range(($colRange -f 165)).Value2

$row167 = $ws.Range(($colRange -f 167)).Value2
$row168 = $ws.Range(($colRange -f 168)).Value2

# ---- write the new pairing back ----
# rows 11 / 12 swap with each other
$ws.Range(($colRange -f 11)).Value2 = $row12
$ws.Range(($colRange -f 12)).Value2 = $row11

# rows 15 / 16 / 17 rotate: new15 = old17, new16 = old15, new17 = old16
$ws.Range(($colRange -f 15)).Value2 = $row17
$ws.Range(($colRange -f 16)).Value2 = $row15
$ws.Range(($colRange -f 17)).Value2 = $row16

# rows 21 / 22 swap with each other
$ws.Range(($colRange -f 21)).Value2 = $row22
$ws.Range(($colRange -f 22)).Value2 = $row21

# rows 101 / 102 swap with each other
$ws.Range(($colRange -f 101)).Value2 = $row102
$ws.Range(($colRange -f 102)).Value2 = $row101

# rows 164 / 165 swap with each other
$ws.Range(($colRange -f 164)).Value2 = $row165
$ws.Range(($colRange -f 165)).Value2 = $row164

# rows 167 / 168 swap with each other
$ws.Range(($colRange -f 167)).Value2 = $row168
$ws.Range(($colRange -f 168)).Value2 = $row167
